# Applies cryptocurrency price/volume updates to Sheet1 (cryptos.xlsx)
# Rows 13 and 14 also swap their Coin name and Link (WrappedEther <-> WrappedliquidstakedEther2.0)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while forcing it to be stored as literal
# text (not auto-converted to a number), and without leaving behind any
# number-format / quote-prefix style change on the cell.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "25.789.86"
$ws.Range("E2").Value = "  +0.13%  "

Set-TextValue $ws.Range("D3") "1.630.02"
$ws.Range("E3").Value = "  -0.28%  "

Set-TextValue $ws.Range("D4") "0.997"
$ws.Range("E4").Value = "  -0.49%  "

Set-TextValue $ws.Range("D5") "213.94"
$ws.Range("E5").Value = "  -0.65%  "

Set-TextValue $ws.Range("D6") "0.502"
$ws.Range("E6").Value = "  -0.37%  "

Set-TextValue $ws.Range("D7") "0.997"
$ws.Range("E7").Value = "  -0.54%  "

Set-TextValue $ws.Range("D8") "0.256"
$ws.Range("E8").Value = "  -0.88%  "

Set-TextValue $ws.Range("D9") "0.0630"
$ws.Range("E9").Value = "  -0.84%  "

Set-TextValue $ws.Range("D10") "19.62"
$ws.Range("E10").Value = "  +0.28%  "

$ws.Range("E11").Value = "  +0.29%  "

Set-TextValue $ws.Range("D12") "4.26"
$ws.Range("E12").Value = "  +0.42%  "

$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue $ws.Range("D13") "1.853.98"
$ws.Range("E13").Value = "  -0.27%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Range("D14") "1.622.90"
$ws.Range("E14").Value = "  -0.72%  "

Set-TextValue $ws.Range("D15") "0.551"
$ws.Range("E15").Value = "  -0.85%  "

$ws.Range("E16").Value = "  -1.08%  "

Set-TextValue $ws.Range("D17") "62.63"
$ws.Range("E17").Value = "  -0.21%  "

Set-TextValue $ws.Range("D18") "25.782.20"
$ws.Range("E18").Value = "  +0.04%  "

Set-TextValue $ws.Range("D19") "0.997"
$ws.Range("E19").Value = "  -0.50%  "

Set-TextValue $ws.Range("D20") "4.43"
$ws.Range("E20").Value = "  -0.07%  "

Set-TextValue $ws.Range("D21") "190.66"
$ws.Range("E21").Value = "  -1.60%  "

Set-TextValue $ws.Range("D22") "9.94"
$ws.Range("E22").Value = "  +0.02%  "

Set-TextValue $ws.Range("D23") "6.29"
$ws.Range("E23").Value = "  -0.02%  "

$ws.Range("E24").Value = "  -0.64%  "

$ws.Range("E25").Value = "  -2.26%  "

Set-TextValue $ws.Range("D26") "142.25"
$ws.Range("E26").Value = "  +1.27%  "

$ws.Range("E27").Value = "  +0.51%  "

Set-TextValue $ws.Range("D28") "6.83"
$ws.Range("E28").Value = "  -0.79%  "

Set-TextValue $ws.Range("D29") "15.46"
$ws.Range("E29").Value = "  -0.29%  "

$ws.Range("E30").Value = "  -0.74%  "

$ws.Range("E31").Value = "  +0.22%  "

$ws.Range("E32").Value = "  -0.49%  "

$ws.Range("E33").Value = "  -0.74%  "

$ws.Range("E34").Value = "  +0.05%  "

$ws.Range("E35").Value = "  +0.11%  "

Set-TextValue $ws.Range("D36") "0.903"
$ws.Range("E36").Value = "  +0.43%  "

Set-TextValue $ws.Range("D37") "1.143.77"
$ws.Range("E37").Value = "  +2.11%  "

$ws.Range("E38").Value = "  -0.75%  "

$ws.Range("E39").Value = "  -1.33%  "

Set-TextValue $ws.Range("D40") "0.0156"
$ws.Range("E40").Value = "  +0.05%  "

Set-TextValue $ws.Range("D41") "0.996"
$ws.Range("E41").Value = "  -0.64%  "

Set-TextValue $ws.Range("D42") "5.61"
$ws.Range("E42").Value = "  +0.60%  "

Set-TextValue $ws.Range("D43") "100.36"
$ws.Range("E43").Value = "  +0.77%  "

Set-TextValue $ws.Range("D44") "0.800"
$ws.Range("E44").Value = "  -0.20%  "

Set-TextValue $ws.Range("D45") "1.764.60"
$ws.Range("E45").Value = "  -0.20%  "

Set-TextValue $ws.Range("D46") "0.0₆0110"
$ws.Range("E46").Value = "  -1.47%  "

Set-TextValue $ws.Range("D47") "55.37"
$ws.Range("E47").Value = "  +0.62%  "

Set-TextValue $ws.Range("D48") "0.0512"
$ws.Range("E48").Value = "  +2.16%  "

$ws.Range("E49").Value = "  +5.29%  "

Set-TextValue $ws.Range("D50") "0.416"
$ws.Range("E50").Value = "  -0.24%  "

Set-TextValue $ws.Range("D51") "7.55"
$ws.Range("E51").Value = "  -0.49%  "
